$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Insert 6 new rows starting at row 12. This pushes the old empty filler
# rows (12-14) down to (18-20) and the old "Goederen..." block (16-25) down
# to (22-31), matching the new layout / dimension (A1:B31).
# ---------------------------------------------------------------------------
$ws.Rows("12:17").Insert()

# Give the freshly inserted rows the same look as their neighbours before
# putting any text in them (copying formats only, so row heights are not
# touched here).
$ws.Range("B3").Copy()
$ws.Range("B12:B15").PasteSpecial(-4122)
$ws.Range("B18").Copy()
$ws.Range("B16:B17").PasteSpecial(-4122)
$ws.Range("B21").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# Row 4 (B4): turn the plain "Ut 07:38-07:57..." note into a rich-text cell
# with a bold header line ("Ut 07:38-07:57") followed by the rest of the
# note in regular weight.
# ---------------------------------------------------------------------------
$b4 = @"
Ut 07:38-07:57
Tabel 20a: incompleet
Tabel 20b:
[n] 6018 >?
[x] 1818
[x] 7318 <?
[x] 9618 >?
[n] 9620(A) >?
"@
$ws.Range("B4").Value = $b4
$ws.Range("B4").Characters(1, 14).Font.Bold = $true
$ws.Range("B4").Characters(15, 1000).Font.Bold = $false

# ---------------------------------------------------------------------------
# New row 12 (B12): bold header "Ut-Amf 07:52-08:07" followed by the
# Tabel 32 note in regular weight.
# ---------------------------------------------------------------------------
$b12 = @"
Ut-Amf 07:52-08:07
Tabel 32a:
[n] 3921
[n] 5723
Tabel 32b:
[x] 5720
[n] 3922
"@
$ws.Range("B12").Value = $b12
$ws.Range("B12").Characters(1, 18).Font.Bold = $true
$ws.Range("B12").Characters(19, 1000).Font.Bold = $false
$ws.Rows("12:12").RowHeight = 99.75

# New row 13 (B13): Tabel 34 note.
$b13 = @"
Tabel 34a:
[x] 3623
[x] 5623
Tabel 34b:
[x] 3618
[x] 5620/7420
"@
$ws.Range("B13").Value = $b13
$ws.Rows("13:13").RowHeight = 85.5

# New row 14 (B14): Tabel 70 note.
$b14 = @"
Tabel 70a:
[x] 3623
[x] 5623
Tabel 70b:
[x] 3618
[x] 5620
[x] 5620/7420
[x] 520
[x] 3620
"@
$ws.Range("B14").Value = $b14
$ws.Rows("14:14").RowHeight = 128.25

# New row 15 (B15): Tabel 80 note.
$b15 = @"
Tabel 80a:
[x] 3623
[x] 5623
Tabel 80b:
[x] 3618
[x] 5620 (1)
[x] 7420
[x] 520
[x] 3620
"@
$ws.Range("B15").Value = $b15
$ws.Rows("15:15").RowHeight = 128.25

# New row 16 (B16): short plain note, formatted like the empty filler rows
# (wrap text, no fill) rather than the highlighted note rows.
$ws.Range("B16").Value = "Amf 08:02-08:17"

# ---------------------------------------------------------------------------
# View state: selection moved to B17.
# ---------------------------------------------------------------------------
$ws.Range("B17").Select()

# ---------------------------------------------------------------------------
# Page setup: paper size / orientation set explicitly (as when printing).
# ---------------------------------------------------------------------------
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
